$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I48").Value = "try solution this morning and some clumpy beads so make new solution and add azide this time"

$ws.Range("I49").Value = "ypos was low this morning so soak micro 25min"

$ws.Range("A50").Value = "D20151105T194435"
$ws.Range("B50").Value = 5
$ws.Range("C50").Value = "V"
$ws.Range("D50").Value = 594.57000000000005
$ws.Range("E50").Value = 4.1223000000000001
$ws.Range("F50").Value = 1198
$ws.Range("G50").Value = 191.46
$ws.Range("H50").Value = 2451

$ws.Range("A51").Value = "D20151105T200655"
$ws.Range("B51").Value = 5
$ws.Range("C51").Value = "V"
$ws.Range("D51").Value = 598.27
$ws.Range("E51").Value = 1.1787000000000001
$ws.Range("F51").Value = 1198
$ws.Range("G51").Value = 195.1
$ws.Range("H51").Value = 2500

$ws.Range("A52").Value = "D20151105T202915"
$ws.Range("B52").Value = 5
$ws.Range("C52").Value = "V"
$ws.Range("D52").Value = 585.6
$ws.Range("E52").Value = 4.1288999999999998
$ws.Range("F52").Value = 1198
$ws.Range("G52").Value = 188
$ws.Range("H52").Value = 2418
$ws.Range("I52").Value = "ypos consistent through files, no real junk, separate beads"

$ws.Range("A53").Value = "D20151105T210041"
$ws.Range("B53").Value = 5
$ws.Range("C53").Value = "H"

$ws.Range("D53").Select()
